# Insert a new summary row at row 2 that computes column averages over the
# (now shifted down) data rows 3:13, labelled "Average" in column A, with the
# rest of the metadata columns (B:G) copied from the first data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift all existing data rows (old 2..12) down to 3..13.
$ws.Rows(2).Insert()

# Label the new summary row.
$ws.Range("A2").Value = "Average"

# Metadata columns B:G are identical across the block, so just copy them
# from the row immediately below (the former row 2, now row 3).
$ws.Range("B2:G2").Value2 = $ws.Range("B3:G3").Value2

# Numeric columns H:CF get a per-column AVERAGE over the 11 data rows
# (3 through 13). Assigning the formula to the whole range lets each cell
# pick up its own relative column reference, same as a fill-right in Excel.
$ws.Range("H2:CF2").Formula = "=AVERAGE(H3:H13)"

# Match the saved selection state from the edit.
$ws.Range("A3").Select()
